# docker network başlıkları düzenlendi
#
# 1) "MULTI-STAGE" paragraph becomes bold (b + bCs), matching the style
#    used by the other section headers in the document (e.g. "Inspection",
#    "Host Mode").
# 2) Several paragraphs that used to be split across two runs -- the main
#    sentence run plus a trailing "(PDF nnn)" run -- get merged back into a
#    single run/string.

$d = $word.ActiveDocument

# --- 1. Bold the "MULTI-STAGE" heading paragraph -------------------------
$rng = $d.Content
$found = $rng.Find.Execute("MULTI-STAGE", $true, $false, $false, $false,
                            $false, $true, 1, $false, "", 0)
if ($found) {
    $para = $rng.Paragraphs(1)
    $full = $para.Range
    $full.Font.Bold = 1
    $full.Font.BoldBi = 1
}

# --- 2. Merge the split "(PDF nnn)" runs back into their preceding run ---

$merges = @(
    "Development Compose file) (PDF 257)",
    "Ilgili rule table a yazmamız lazım. (PDF 326)",
    "SERVICE DISCOVERY (PDF 332)",
    "DOCKER COMPOSE (PDF 456)"
)

foreach ($needle in $merges) {
    $d.Content.Find.Execute($needle, $false, $false, $false, $false, $false,
                             $true, 1, $false, $needle, 2) | Out-Null
}

# --- 3. Add suppressAutoHyphens to the Normal style's paragraph format ---
$normal = $d.Styles.Item(1)
$normal.ParagraphFormat.Hyphenation = 0
